$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = '(''باسلْ ليسك – الحجر واللسان'', [''{4}{G}{G}{G}'', ''مخلوق – باسلْ ليسك'', ''أينما يقوم باسلْ ليسك — الحجر واللسان، بتوجيه ضرر قتالي لمخلوق اَخر، قم بتحطيم هذا المخلوق المضروب عند نهاية القتال.'', ''الحد — على جميع المخلوقات القادرة على سد هجوم باسلْ ليسك — الحجر واللسان، عمل ذلك. (لديك الحد طالما وجد سبع بطاقات أو أكثر في مقبرتك).'', ''4/5''])'

$ws.Range("A3:A7").EntireRow.Delete()
